# Refresh the cryptos list (Price / Volume(1h) columns) to match the
# source feed snapshot taken on Wed Nov 20 07:31:21 UTC 2024.
#
# All Price/Volume cells are stored as plain text in this sheet (values
# such as "1.00", "3.107.77" or "  +0.94%  " must keep their exact
# characters -- trailing zeros, thousand-dot grouping, padding spaces,
# etc). Force General/Price/Volume columns to Text format before writing
# so Excel does not silently coerce them to numbers and strip formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold text-formatted numbers/percentages.
$ws.Columns("D:E").NumberFormat = "@"

$ws.Range('D2').Value = '92.736.95'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '3.110.80'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '235.47'
$ws.Range('E5').Value = '  -3.12%  '
$ws.Range('D6').Value = '613.49'
$ws.Range('E6').Value = '  -0.69%  '
$ws.Range('E7').Value = '  -2.30%  '
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '3.107.33'
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('D11').Value = '0.784'
$ws.Range('E11').Value = '  +3.71%  '
$ws.Range('E12').Value = '  -3.67%  '
$ws.Range('E13').Value = '  -3.70%  '
$ws.Range('D14').Value = '92.491.04'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').Value = '33.89'
$ws.Range('E15').Value = '  -4.21%  '
$ws.Range('E16').Value = '  -3.59%  '
$ws.Range('D18').Value = '3.116.39'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('D19').Value = '3.80'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').Value = '14.51'
$ws.Range('E20').Value = '  -3.26%  '
$ws.Range('D21').Value = '5.83'
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('D23').Value = '439.21'
$ws.Range('E23').Value = '  -3.84%  '
$ws.Range('D24').Value = '9.10'
$ws.Range('E24').Value = '  -1.13%  '
$ws.Range('D25').Value = '8.21'
$ws.Range('E25').Value = '  +5.22%  '
$ws.Range('D26').Value = '5.58'
$ws.Range('E26').Value = '  -6.37%  '
$ws.Range('D27').Value = '85.38'
$ws.Range('E27').Value = '  -4.93%  '
$ws.Range('E28').Value = '  -1.90%  '
$ws.Range('D29').Value = '3.274.68'
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('E31').Value = '  +8.03%  '
$ws.Range('D32').Value = '0.235'
$ws.Range('E32').Value = '  +2.45%  '
$ws.Range('D33').Value = '0.126'
$ws.Range('E33').Value = '  -12.52%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value = '1.03'
$ws.Range('E34').Value = '  -30.38%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '9.15'
$ws.Range('E35').Value = '  -2.40%  '
$ws.Range('D36').Value = '8.06'
$ws.Range('E36').Value = '  +7.19%  '
$ws.Range('D37').Value = '0.164'
$ws.Range('E37').Value = '  -5.08%  '
$ws.Range('D38').Value = '25.68'
$ws.Range('E38').Value = '  -2.89%  '
$ws.Range('D39').Value = '4.00'
$ws.Range('E39').Value = '  +3.46%  '
$ws.Range('D40').Value = '1.90'
$ws.Range('E40').Value = '  -5.19%  '
$ws.Range('D41').Value = '23.91'
$ws.Range('E41').Value = '  +7.70%  '
$ws.Range('E42').Value = '  -2.73%  '
$ws.Range('D43').Value = '463.89'
$ws.Range('E43').Value = '  -6.09%  '
$ws.Range('D44').Value = '0.427'
$ws.Range('E44').Value = '  -2.76%  '
$ws.Range('E45').Value = '  -1.91%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').Value = '159.55'
$ws.Range('E47').Value = '  +1.63%  '
$ws.Range('D48').Value = '0.682'
$ws.Range('E48').Value = '  -3.83%  '
$ws.Range('E49').Value = '  -5.12%  '
$ws.Range('E50').Value = '  -2.95%  '
$ws.Range('D51').Value = '43.80'
$ws.Range('E51').Value = '  -0.55%  '
